$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value that was bumped by one
# day (45180 -> 45181) for every data row (rows 2 through 300).
$range = $ws.Range("C2:C300")
$range.Value = 45181
